# Update the "Notes" sheet: correct the description/source text, add a
# source-link row, update the licensing text, and add a licensing-info
# link row. This shifts the trailing notes (email/copyright) down.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Notes")

# --- Description (row 2) ---
$ws.Range("A2").Value = "Description: Poverty Headcount"

# --- Source (row 4) ---
$ws.Range("A4").Value = "Source: Profiles of higher local governments 2014 - Uganda Bureau of Statistics"

# --- Insert a new row 5 for the source link (pushes the blank line and
#     everything below it down by one row) ---
$ws.Rows.Item(5).Insert()
$ws.Range("A5").Value = "Source-link: http://www.ubos.org/onlinefiles/uploads/ubos/2009_HLG_%20Abstract_printed/CIS+UPLOADS/Profiles%20of%20Higher%20Local%20Governments_June_2014.pdf"

# After the insert above:
#   row 6  = blank
#   row 7  = "Notes:"
#   row 8  = blank
#   row 9  = Data-wide-value note
#   row 10 = blank
#   row 11 = blank
#   row 12 = blank
#   row 13 = Datahub download note
#   row 14 = old license text
#   row 15 = email contact
#   row 16 = copyright

# --- Update the license text in place (row 14) ---
$ws.Range("A14").Value = "It is licensed under a Creative Commons Attribution 4.0 International license."

# --- Insert a new row 15 for the licensing info link (pushes email and
#     copyright down by one more row) ---
$ws.Rows.Item(15).Insert()
$ws.Range("A15").Value = "More information on licensing is available here: https://creativecommons.org/licenses/by/4.0/"

# Final layout:
#   row 16 = email contact (unchanged content, shifted)
#   row 17 = copyright (unchanged content, shifted)
